$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name / Link columns (plain text, swapped rows)
$bcUpdates = @(
    @("B29", "PEPE"),
    @("C29", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"),
    @("B30", "Monero"),
    @("C30", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("B42", "Bittensor"),
    @("C42", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"),
    @("B43", "Aave"),
    @("C43", "https://coinranking.com/coin/ixgUfzmLR+aave-aave"),
    @("B44", "Stellar"),
    @("C44", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"),
    @("B45", "Hedera"),
    @("C45", "https://coinranking.com/coin/jad286TjB+hedera-hbar"),
    @("B46", "Polygon"),
    @("C46", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic")
)
foreach ($pair in $bcUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Update Price / Volume(1h) columns. Force text so numeric-looking
# strings (e.g. "134.04") are not auto-converted to numbers, matching
# the source workbook where these are inline/shared strings, then
# restore the default "Normal" style so no stray per-cell format sticks.
$deUpdates = @(
    @("D2", "58.460.48"),
    @("E2", "  +1.21%  "),
    @("D3", "2.370.12"),
    @("E3", "  +1.56%  "),
    @("E4", "  +0.04%  "),
    @("D5", "549.46"),
    @("E5", "  +1.06%  "),
    @("D6", "134.04"),
    @("E6", "  -0.93%  "),
    @("E7", "  -0.04%  "),
    @("E8", "  +1.22%  "),
    @("D9", "0.108"),
    @("E9", "  +5.76%  "),
    @("D10", "5.70"),
    @("E10", "  +4.50%  "),
    @("E11", "  -1.19%  "),
    @("E12", "  -0.43%  "),
    @("D13", "24.28"),
    @("E13", "  +2.26%  "),
    @("D14", "2.792.84"),
    @("E14", "  +1.65%  "),
    @("D15", "58.397.20"),
    @("E15", "  +1.18%  "),
    @("E16", "  +3.07%  "),
    @("D17", "2.352.23"),
    @("E17", "  +0.62%  "),
    @("D18", "11.13"),
    @("E18", "  +4.51%  "),
    @("E19", "  +2.15%  "),
    @("D20", "333.00"),
    @("E20", "  -0.43%  "),
    @("D21", "7.05"),
    @("E21", "  +5.50%  "),
    @("E22", "  +0.27%  "),
    @("D23", "64.03"),
    @("E23", "  +3.23%  "),
    @("D24", "0.169"),
    @("E24", "  +1.15%  "),
    @("D25", "0.998"),
    @("E25", "  -0.32%  "),
    @("D26", "8.34"),
    @("E26", "  -1.97%  "),
    @("E27", "  -5.03%  "),
    @("E28", "  +0.53%  "),
    @("D29", "0.0₃0748"),
    @("E29", "  +2.14%  "),
    @("D30", "170.34"),
    @("E30", "  -0.04%  "),
    @("D31", "6.19"),
    @("E31", "  +0.78%  "),
    @("D32", "18.49"),
    @("E32", "  +0.16%  "),
    @("E33", "  -2.27%  "),
    @("E35", "  +0.20%  "),
    @("D36", "4.19"),
    @("E36", "  +0.23%  "),
    @("D37", "1.25"),
    @("E37", "  -0.55%  "),
    @("D38", "40.44"),
    @("E38", "  +3.36%  "),
    @("D39", "0.424"),
    @("E39", "  +13.25%  "),
    @("E40", "  -1.03%  "),
    @("D41", "3.72"),
    @("E41", "  +2.55%  "),
    @("D42", "290.90"),
    @("E42", "  +1.60%  "),
    @("D43", "141.34"),
    @("E43", "  -3.10%  "),
    @("D44", "0.0961"),
    @("E44", "  +2.57%  "),
    @("D45", "0.0516"),
    @("E45", "  +2.55%  "),
    @("D46", "0.415"),
    @("E46", "  +9.00%  "),
    @("D47", "0.567"),
    @("E47", "  +1.20%  "),
    @("D48", "18.73"),
    @("E48", "  -1.93%  "),
    @("D49", "0.0224"),
    @("E49", "  +3.30%  "),
    @("D50", "11.04"),
    @("E50", "  -0.19%  ")
)
foreach ($pair in $deUpdates) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.Style = "Normal"
}
